# Refresh the crypto price table with the latest scraped values.
# Price (D) and Volume(1h) (E) columns are stored as literal text in the
# sheet, so a leading apostrophe is prepended to stop Excel from
# reinterpreting the numeric-looking / percentage-looking strings as
# numbers (same as the "quote prefix" you get by typing '257.96 by hand).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'257.96"
$ws.Range('E2').Value = "'0.05%"
$ws.Range('D3').Value = "'27.00"
$ws.Range('E3').Value = "'-0.55%"
$ws.Range('D4').Value = "'4.658"
$ws.Range('E4').Value = "'-4.42%"
$ws.Range('D6').Value = "'6.645"
$ws.Range('E6').Value = "'-0.63%"
$ws.Range('D7').Value = "'0.8552"
$ws.Range('E7').Value = "'-1.56%"
$ws.Range('D8').Value = "'0.9474"
$ws.Range('E8').Value = "'-1.24%"
$ws.Range('D9').Value = "'0.1403"
$ws.Range('E9').Value = "'-0.59%"
$ws.Range('D10').Value = "'0.05148"
$ws.Range('E10').Value = "'45.37%"
$ws.Range('D11').Value = "'0.07092"
$ws.Range('E11').Value = "'-1.25%"
$ws.Range('D12').Value = "'0.03106"
$ws.Range('E12').Value = "'-1.44%"
$ws.Range('D13').Value = "'0.09143"
$ws.Range('E13').Value = "'-1.07%"
$ws.Range('D14').Value = "'0.001528"
$ws.Range('E14').Value = "'-0.82%"
$ws.Range('D15').Value = "'0.0006068"
$ws.Range('E15').Value = "'0.44%"
$ws.Range('D16').Value = "'0.006176"
$ws.Range('E16').Value = "'3.21%"
$ws.Range('D17').Value = "'3.502"
$ws.Range('E17').Value = "'0.49%"
$ws.Range('E18').Value = "'-2.20%"
$ws.Range('D20').Value = "'0.3054"
$ws.Range('E20').Value = "'-2.93%"
$ws.Range('D21').Value = "'0.1278"
$ws.Range('E21').Value = "'-2.24%"
$ws.Range('D22').Value = "'3.815"
$ws.Range('E22').Value = "'7.90%"
$ws.Range('E23').Value = "'-0.36%"
$ws.Range('D24').Value = "'0.001222"
$ws.Range('E24').Value = "'-0.12%"
$ws.Range('D25').Value = "'0.004298"
$ws.Range('E25').Value = "'-4.85%"
$ws.Range('D27').Value = "'0.0001936"
$ws.Range('E27').Value = "'29.82%"
$ws.Range('D40').Value = "'0.03830"
$ws.Range('E40').Value = "'-0.06%"
$ws.Range('B41').Value = "BKEXToken"
$ws.Range('C41').Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range('D41').Value = "'0.1103"
$ws.Range('E41').Value = "'0.06%"
$ws.Range('B42').Value = "KickToken"
$ws.Range('C42').Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range('D42').Value = "'0.006259"
$ws.Range('E42').Value = "'-4.89%"
$ws.Range('D43').Value = "'0.002370"
$ws.Range('E43').Value = "'7.69%"
$ws.Range('D44').Value = "'0.01396"
$ws.Range('E44').Value = "'32.77%"
$ws.Range('D45').Value = "'0.00005395"
$ws.Range('E45').Value = "'-1.83%"
$ws.Range('E46').Value = "'-0.05%"
$ws.Range('D47').Value = "'0.05098"
$ws.Range('E47').Value = "'-53.28%"
$ws.Range('D48').Value = "'0.2530"
$ws.Range('E48').Value = "'11,788.49%"
$ws.Range('D49').Value = "'0.00002099"
$ws.Range('E49').Value = "'-0.05%"
$ws.Range('D50').Value = "'0.0001999"
$ws.Range('E50').Value = "'-0.05%"
